$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update date in A1
$ws.Range("A1").Value = 45436

# Update price column (D) values
$ws.Range("D19").Value = 2296.529
$ws.Range("D20").Value = 2515.757
$ws.Range("D21").Value = 2604.87
$ws.Range("D22").Value = 2714.087
$ws.Range("D23").Value = 2746.031
$ws.Range("D24").Value = 3581.307
$ws.Range("D25").Value = 4034.989
$ws.Range("D26").Value = 4809.861
$ws.Range("D27").Value = 5407.288
$ws.Range("D28").Value = 5885.863
$ws.Range("D29").Value = 6431.885
$ws.Range("D30").Value = 7130.494
$ws.Range("D31").Value = 7708.642
$ws.Range("D32").Value = 8078.004
$ws.Range("D33").Value = 9178.103999999999
$ws.Range("D34").Value = 2874.381
$ws.Range("D35").Value = 10101.518
$ws.Range("D36").Value = 11322.051
$ws.Range("D37").Value = 12125.037
$ws.Range("D43").Value = 3666.504
$ws.Range("D44").Value = 3942.648
$ws.Range("D45").Value = 4008.496
$ws.Range("D46").Value = 4255.802
$ws.Range("D47").Value = 4665.322
$ws.Range("D48").Value = 5219.386
$ws.Range("D49").Value = 6247.204
$ws.Range("D50").Value = 7266.987
$ws.Range("D51").Value = 8078.004
$ws.Range("D52").Value = 8905.083000000001
$ws.Range("D53").Value = 9916.833000000001
$ws.Range("D54").Value = 11474.618
$ws.Range("D55").Value = 12205.336
$ws.Range("D56").Value = 13771.148
$ws.Range("D57").Value = 15409.231
$ws.Range("D58").Value = 17617.436
$ws.Range("D59").Value = 19279.613
$ws.Range("D66").Value = 5386.553
$ws.Range("D67").Value = 5434.414
$ws.Range("D68").Value = 5885.863
$ws.Range("D69").Value = 6247.204
$ws.Range("D70").Value = 6704.904
$ws.Range("D71").Value = 8045.883
$ws.Range("D72").Value = 8816.75
$ws.Range("D73").Value = 11105.249
$ws.Range("D74").Value = 11675.354
$ws.Range("D75").Value = 13313.453
$ws.Range("D76").Value = 14397.484
$ws.Range("D77").Value = 16236.309
$ws.Range("D78").Value = 17513.047
$ws.Range("D79").Value = 19279.613
$ws.Range("D80").Value = 20371.665
$ws.Range("D81").Value = 22965.306
$ws.Range("D82").Value = 25141.369
$ws.Range("D83").Value = 26980.216
$ws.Range("D90").Value = 17812.872
$ws.Range("D91").Value = 40871.806
